$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Patient Id"
$ws.Range("B1").Value = "Trial"
$ws.Range("C1").Value = "Correct"
$ws.Range("D1").Value = "Elapsed Time"
$ws.Range("E1").Value = "Date"

# --- Data rows ---
# Columns: A = Patient Id (number), B = Trial (number),
#          C = Correct (text "false"), D = Elapsed Time (number),
#          E = Date (text "2019-10-21")
$patientIds = @(4, 4, 1, 3, 3, 1, 2, 2, 0, 0)
$trials     = @(1, 2, 1, 1, 2, 2, 1, 2, 1, 2)

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $patientIds[$i]
    $ws.Range("B$r").Value = $trials[$i]
    $ws.Range("C$r").Value = "'false"
    $ws.Range("D$r").Value = 92
    $ws.Range("E$r").Value = "'2019-10-21"
}

# Strip the quote-prefix styling the apostrophe trick above applies, so the
# cells stay plain text (t="s") without carrying a quotePrefix cell style.
$ws.Range("C2:C11").ClearFormats()
$ws.Range("E2:E11").ClearFormats()
